{"js": "// The document contains four \"tc\" id fields of the form:\n//   <id>p084r_aN</id>\n// split across three runs (the \"<id>\" / \"aN-body\" / \"</id>\" runs), where\n// the middle run carries the bare id value (e.g. \"p084r_a1\"). The newly\n// downloaded tc/tcn/tl content renumbers these ids by dropping the \"a\"\n// prefix letter (p084r_a1 -> p084r_1, p084r_a3 -> p084r_3, etc.), and the\n// three runs collapse into the single surviving \"<id>...</id>\" run.\n//\n// We search for the whole \"<id>p084r_aN</id>\" text (Word's search can span\n// multiple runs) and replace the matched range's text in one shot, which\n// merges/collapses the three runs into one, matching the target edit.\nconst idMap = {\n  \"p084r_a1\": \"p084r_1\",\n  \"p084r_a3\": \"p084r_3\",\n  \"p084r_a4\": \"p084r_4\",\n  \"p084r_a6\": \"p084r_6\",\n};\n\nfor (const oldId in idMap) {\n  const newId = idMap[oldId];\n  const searchResults = context.document.body.search(`<id>${oldId}</id>`, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  searchResults.load(\"items\");\n  await context.sync();\n\n  for (const range of searchResults.items) {\n    range.insertText(`<id>${newId}</id>`, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# The document contains four \"tc\" id fields of the form:\n#   <id>p084r_aN</id>\n# split across three runs (the \"<id>\" / \"aN-body\" / \"</id>\" runs), where\n# the middle run carries the bare id value (e.g. \"p084r_a1\"). The newly\n# downloaded tc/tcn/tl content renumbers these ids by dropping the \"a\"\n# prefix letter (p084r_a1 -> p084r_1, p084r_a3 -> p084r_3, etc.), and the\n# three runs collapse into the single surviving \"<id>...</id>\" run.\n#\n# For each id we Find the whole \"<id>p084r_aN</id>\" text (Find spans runs),\n# then set the matched Range's Text directly, which merges/collapses the\n# three runs into one, matching the target edit.\n\n$d = $word.ActiveDocument\n\n$idMap = [ordered]@{\n    \"p084r_a1\" = \"p084r_1\"\n    \"p084r_a3\" = \"p084r_3\"\n    \"p084r_a4\" = \"p084r_4\"\n    \"p084r_a6\" = \"p084r_6\"\n}\n\nforeach ($oldId in $idMap.Keys) {\n    $newId = $idMap[$oldId]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = \"<id>$oldId</id>\"\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = \"\"\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $found = $find.Execute()\n    if ($found) {\n        $rng.Text = \"<id>$newId</id>\"\n    }\n}\n"}
